{"js": "// Office.js (Word JavaScript API) edit script.\n// Body of: async (context) => { ... }\n//\n// Goal (from the diff): in five specific bullet/impact paragraphs, wrap the\n// quantitative metrics (percentages, dollar amounts, large numbers) in their\n// own runs with bold + color (#2C3E50) formatting, while leaving the rest of\n// the paragraph's text as plain runs. Two of the paragraphs share the same\n// \"Achieved 87% ... 71%\" prefix so we disambiguate by exact full text.\n\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items/text\");\nawait context.sync();\n\n// Exact full paragraph texts (as they exist BEFORE this edit) mapped to the\n// ordered list of metric substrings that must become bold + colored.\nconst targets = [\n  {\n    text:\n      \"\\u2022 Discovered systematic race coding errors affecting all Black and Asian-American voters, developed geospatial machine learning algorithms improving demographic classification accuracy from 23% to 64%\",\n    metrics: [\"23%\", \"64%\"],\n  },\n  {\n    text:\n      \"\\u2022 Achieved 87% prediction accuracy for voter turnout vs. industry standard of 71%, reducing polling error margins from \\u00B14.2% to \\u00B12.1%\",\n    metrics: [\"87%\", \"71%\", \"\\u00B14.2%\", \"\\u00B12.1%\"],\n  },\n  {\n    text: \"\\u2022 Wrote RFP and analyzed bids from 1,200 vendors for research platform development\",\n    metrics: [\"1,200\"],\n  },\n  {\n    text:\n      \"\\u2022 Created comprehensive meta-analysis framework handling millions of survey responses that became the $400M Polling Consortium Database at The Analyst Institute, now valued at $1B+\",\n    metrics: [\"$400M\", \"$1B\"],\n  },\n  {\n    text: \"\\u2022 Algorithm reduced mapping costs by 73.5%, saving campaigns and organizations $4.7M\",\n    metrics: [\"73.5%\", \"$4.7M\"],\n  },\n  {\n    text: \"\\u2022 Achieved 87% prediction accuracy for voter turnout vs. industry standard of 71%\",\n    metrics: [\"87%\", \"71%\"],\n  },\n];\n\n// Find the matching paragraph Word Proxy objects up front (by exact text),\n// in document order, matching each target to the first not-yet-used\n// paragraph with identical text (handles the duplicate-prefix case because\n// the long variant's text != the short variant's text).\nconst used = new Set();\nconst matchedParagraphs = [];\nfor (const target of targets) {\n  let foundIndex = -1;\n  for (let i = 0; i < paragraphs.items.length; i++) {\n    if (used.has(i)) continue;\n    if (paragraphs.items[i].text === target.text) {\n      foundIndex = i;\n      break;\n    }\n  }\n  if (foundIndex === -1) {\n    throw new Error(\"Could not locate paragraph for text: \" + target.text);\n  }\n  used.add(foundIndex);\n  matchedParagraphs.push(paragraphs.items[foundIndex]);\n}\n\n// Apply bold + color to each metric substring, scoped to its own paragraph\n// (paragraph.search() only searches within that paragraph's range, so the\n// unrelated occurrences of the same percentages elsewhere in the doc are\n// left untouched).\nfor (let t = 0; t < targets.length; t++) {\n  const paragraph = matchedParagraphs[t];\n  const metrics = targets[t].metrics;\n  for (const metric of metrics) {\n    const results = paragraph.search(metric, { matchCase: true });\n    results.load(\"items\");\n    await context.sync();\n    for (const result of results.items) {\n      result.font.bold = true;\n      result.font.color = \"#2C3E50\";\n    }\n    await context.sync();\n  }\n}\n", "ps1": "# Word COM interop (PowerShell-style) edit script.\n# $word / $doc / $app resolve; the document is open as $word.ActiveDocument.\n#\n# Goal (from the diff): in five specific bullet/impact paragraphs, wrap the\n# quantitative metrics (percentages, dollar amounts, large numbers) in their\n# own runs with bold + color (#2C3E50) formatting, while leaving the rest of\n# the paragraph's text as plain runs. Two of the paragraphs share the same\n# \"Achieved 87% ... 71%\" prefix so we disambiguate by exact full text.\n\nfunction Get-BgrColor($hex) {\n    $r = [Convert]::ToInt32($hex.Substring(0, 2), 16)\n    $g = [Convert]::ToInt32($hex.Substring(2, 2), 16)\n    $b = [Convert]::ToInt32($hex.Substring(4, 2), 16)\n    return $b * 65536 + $g * 256 + $r\n}\n\n$highlightColor = Get-BgrColor \"2C3E50\"\n\n# Bold + color the first occurrence of $text found inside $paragraph's own\n# range (Find.Execute narrows the range to the hit, so formatting only\n# touches that substring's run).\nfunction Format-Metric($paragraph, $text) {\n    $rng = $paragraph.Range\n    $find = $rng.Find\n    $find.ClearFormatting()\n    $find.Text = $text\n    $find.MatchCase = $true\n    $find.MatchWildcards = $false\n    $find.Forward = $true\n    $find.Wrap = 0\n    $found = $find.Execute()\n    if ($found) {\n        $rng.Font.Bold = $true\n        $rng.Font.Color = $highlightColor\n    }\n    return $found\n}\n\n$d = $word.ActiveDocument\n\n# Exact full paragraph texts (as they exist BEFORE this edit) mapped to the\n# ordered list of metric substrings that must become bold + colored.\n$targets = @(\n    @{\n        Text    = \"\u2022 Discovered systematic race coding errors affecting all Black and Asian-American voters, developed geospatial machine learning algorithms improving demographic classification accuracy from 23% to 64%\"\n        Metrics = @(\"23%\", \"64%\")\n    },\n    @{\n        Text    = \"\u2022 Achieved 87% prediction accuracy for voter turnout vs. industry standard of 71%, reducing polling error margins from \u00b14.2% to \u00b12.1%\"\n        Metrics = @(\"87%\", \"71%\", \"\u00b14.2%\", \"\u00b12.1%\")\n    },\n    @{\n        Text    = \"\u2022 Wrote RFP and analyzed bids from 1,200 vendors for research platform development\"\n        Metrics = @(\"1,200\")\n    },\n    @{\n        Text    = \"\u2022 Created comprehensive meta-analysis framework handling millions of survey responses that became the `$400M Polling Consortium Database at The Analyst Institute, now valued at `$1B+\"\n        Metrics = @('$400M', '$1B')\n    },\n    @{\n        Text    = \"\u2022 Algorithm reduced mapping costs by 73.5%, saving campaigns and organizations `$4.7M\"\n        Metrics = @(\"73.5%\", '$4.7M')\n    },\n    @{\n        Text    = \"\u2022 Achieved 87% prediction accuracy for voter turnout vs. industry standard of 71%\"\n        Metrics = @(\"87%\", \"71%\")\n    }\n)\n\n# Snapshot paragraph text up front (by 1-based COM index) so lookups below\n# are against the ORIGINAL text, not text that earlier edits have already\n# split into multiple runs (splitting runs does not change Paragraphs.Count\n# or paragraph ordinals, but re-reading .Range.Text after edits is extra\n# work we don't need).\n$paraCount = $d.Paragraphs.Count\n$paraText = @{}\nfor ($i = 1; $i -le $paraCount; $i++) {\n    $paraText[$i] = $d.Paragraphs.Item($i).Range.Text.TrimEnd([char]13)\n}\n\n$usedIndexes = @{}\n\nforeach ($target in $targets) {\n    $foundIndex = -1\n    for ($i = 1; $i -le $paraCount; $i++) {\n        if ($usedIndexes.ContainsKey($i)) { continue }\n        if ($paraText[$i] -eq $target.Text) {\n            $foundIndex = $i\n            break\n        }\n    }\n    if ($foundIndex -eq -1) {\n        throw \"Could not locate paragraph for text: $($target.Text)\"\n    }\n    $usedIndexes[$foundIndex] = $true\n\n    $paragraph = $d.Paragraphs.Item($foundIndex)\n    foreach ($metric in $target.Metrics) {\n        Format-Metric $paragraph $metric | Out-Null\n    }\n}\n\nWrite-Output \"done\"\n"}
